$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.26600003242492676,
    0.13200020790100098,
    0.06400012969970703,
    0.05500006675720215,
    0.048999786376953125,
    0.05900001525878906,
    0.05200004577636719,
    0.05900001525878906,
    0.06299996376037598,
    0.05200004577636719,
    0.06100010871887207,
    0.10099983215332031,
    0.06599998474121094,
    0.059999942779541016,
    0.06200003623962402,
    0.06600022315979004,
    0.05299997329711914,
    0.06500005722045898,
    0.059999942779541016,
    0.054000139236450195,
    0.06399989128112793,
    0.06200003623962402,
    0.05800008773803711,
    0.06500005722045898,
    0.11299991607666016,
    0.06599998474121094,
    0.06000018119812012,
    0.059999942779541016,
    0.059999942779541016,
    0.0559999942779541,
    0.06200003623962402,
    0.06299996376037598,
    0.07500004768371582,
    0.06299996376037598,
    0.06000018119812012,
    0.06599998474121094,
    0.054000139236450195,
    0.06500005722045898,
    0.06099987030029297,
    0.06399989128112793,
    0.0559999942779541,
    0.06000018119812012,
    0.07400012016296387,
    0.06599998474121094,
    0.06700015068054199,
    0.05499982833862305,
    0.06699991226196289,
    0.06599998474121094,
    0.06099987030029297,
    0.06700015068054199,
    0.07200002670288086,
    0.06400012969970703,
    0.056999921798706055,
    0.06399989128112793,
    0.06399989128112793,
    0.06000018119812012,
    0.07400012016296387,
    0.07500004768371582,
    0.06500005722045898,
    0.07699990272521973,
    0.07500004768371582,
    0.09800004959106445,
    0.06799983978271484,
    0.07700014114379883,
    0.0690000057220459,
    0.07999992370605469,
    0.09200000762939453,
    0.09200000762939453,
    0.07899999618530273,
    0.07299995422363281,
    0.07699990272521973,
    0.08099985122680664,
    0.08499979972839355,
    0.08000016212463379,
    0.07400012016296387,
    0.0820000171661377,
    0.11400008201599121,
    0.08499979972839355,
    0.0819997787475586,
    0.07599997520446777,
    0.07799983024597168,
    0.07800006866455078,
    0.07200002670288086,
    0.07200002670288086,
    0.13699984550476074,
    0.07899999618530273,
    0.08500003814697266,
    0.09099984169006348,
    0.08299994468688965,
    0.08299994468688965,
    0.13200020790100098,
    0.07699990272521973,
    0.0839998722076416,
    0.07999992370605469,
    0.08800005912780762,
    0.07299995422363281,
    0.07400012016296387,
    0.09999990463256836,
    0.08500003814697266,
    0.0839998722076416,
    0.07499980926513672,
    0.09200000762939453,
    0.10299992561340332,
    0.0969998836517334,
    0.08000016212463379,
    0.0820000171661377,
    0.07699990272521973,
    0.10199999809265137,
    0.07699990272521973,
    0.07399988174438477,
    0.07999992370605469,
    0.09100008010864258,
    0.0820000171661377,
    0.07500004768371582,
    0.07500004768371582,
    0.07999992370605469,
    0.07599997520446777,
    0.07599997520446777,
    0.08499979972839355,
    0.09099984169006348
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}
